# Insert two new columns before column C (old C shifts to E), mirroring the
# addition of two new rating-period columns ("Jun_15" and "Jun_17") ahead of
# the existing "Jun_13"/"Jun_10" columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1:D1").EntireColumn.Insert()

# Re-apply the original column width (raw OOXML width ~8.0) to the two new
# columns as well as the shifted-over former column C, so all three keep a
# consistent custom width.
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 7.1

# The new columns start out blank after the insert; fill the data rows with
# the same "UN" placeholder used throughout the rest of the table.
$ws.Range("C2:D27").Value2 = "UN"

# Header row: shift old header text right two slots and add the two new
# period headers.
$ws.Range("D1").Value2 = "Jun_13"
$ws.Range("C1").Value2 = "Jun_15"
$ws.Range("B1").Value2 = "Jun_17"
